$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1 / first sheet) - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 542
$ws1.Range("F4").Value = 194
$ws1.Range("F6").Value = 500
$ws1.Range("F7").Value = 100
$ws1.Range("F8").Value = 115
$ws1.Range("F9").Value = 42
$ws1.Range("F10").Value = 6670
$ws1.Range("F12").Value = 363
$ws1.Range("F13").Value = 2906
$ws1.Range("F14").Value = 188
$ws1.Range("F15").Value = 324
$ws1.Range("F17").Value = 533

# Sheet "全部类型" (sheetId 4 / fourth sheet) - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 542
$ws4.Range("F6").Value = 194
$ws4.Range("F8").Value = 500
$ws4.Range("F9").Value = 100
$ws4.Range("F10").Value = 115
$ws4.Range("F11").Value = 42
$ws4.Range("F13").Value = 6670
$ws4.Range("F16").Value = 363
$ws4.Range("F17").Value = 2906
$ws4.Range("F18").Value = 188
$ws4.Range("F19").Value = 324
$ws4.Range("F21").Value = 533
